# Update the cached "Date and time" auto-field text (datetimeFigureOut
# placeholder) in the slide master and every slide layout:
# 2022/2/19 -> 2022/2/27.
$p = $ppt.ActivePresentation
$m = $p.SlideMaster

for ($i = 1; $i -le $m.Shapes.Count; $i++) {
    $sh = $m.Shapes.Item($i)
    if ($sh.HasTextFrame) {
        if ($sh.PlaceholderFormat.Type -eq 16) {
            $sh.TextFrame.TextRange.Text = "2022/2/27"
        }
    }
}

for ($j = 1; $j -le $m.CustomLayouts.Count; $j++) {
    $lay = $m.CustomLayouts.Item($j)
    for ($i = 1; $i -le $lay.Shapes.Count; $i++) {
        $sh = $lay.Shapes.Item($i)
        if ($sh.HasTextFrame) {
            if ($sh.PlaceholderFormat.Type -eq 16) {
                $sh.TextFrame.TextRange.Text = "2022/2/27"
            }
        }
    }
}

# Slide 1: re-center and shift the "EG STORE" textbox to the right, and
# remove the plain grey sidebar rectangle that used to sit behind it.
$s = $p.Slides.Item(1)

$egShape = $s.Shapes.Item("文本框 3")
$egShape.Left = 323.40002
$egShape.TextFrame.TextRange.ParagraphFormat.Alignment = 2

$rectShape = $s.Shapes.Item("矩形 4")
$rectShape.Delete()
